# Commit: "update data add server"
#
# 1. Add a new expense row (row 13) to the "Budget" sheet: a date, a
#    description ("an kem " - a new shared string), and an amount (20),
#    which flows into the existing SUM formula in row 14.
# 2. Make "Budget" the active / selected sheet (it was "Sheet1" before).

$wb = $excel.ActiveWorkbook
$budget = $wb.Worksheets.Item("Budget")

# --- Add the new row of data on the Budget sheet -------------------------
# Copy the date formatting from an existing date cell in column B so the
# new cell reuses the existing date style instead of creating a new one.
$budget.Range("B9").Copy()
$budget.Range("B13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$budget.Range("B13").Value = 42454          # 3/25/2016
$budget.Range("C13").Value = "ăn kem "
$budget.Range("E13").Value = 20

# --- Switch the active sheet / selection to Budget ------------------------
$budget.Activate()
$budget.Range("E14").Select()
